# Swap the "Steps" / "Expected Results" content between the TC2 block
# (row 18: detalhar diária) and the TC4 block (row 32: analisar prestação
# de contas), so that TC2 now shows the "analisar prestação de contas"
# content and TC4 now shows the "detalhar diária" content. The TC3 block
# and the Test Case IDs themselves (TC2/TC3/TC4 labels) remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tc2Steps    = $ws.Range("B18").Value2
$tc2Expected = $ws.Range("D18").Value2
$tc4Steps    = $ws.Range("B32").Value2
$tc4Expected = $ws.Range("D32").Value2

$ws.Range("B18").Value = $tc4Steps
$ws.Range("D18").Value = $tc4Expected
$ws.Range("B32").Value = $tc2Steps
$ws.Range("D32").Value = $tc2Expected
